# Update investment-cost results with newly-computed server values.
# Sheet mapping: Worksheets.Item(1) = "2025", Item(2) = "2030", Item(3) = "2035",
#                Item(4) = "2040", Item(5) = "2045", Item(6) = "2050" (unchanged)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "2025"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item(1)
$ws.Range("B2").Value = 10372.65132737054
$ws.Range("E2").Value = 289260.5393052954
$ws.Range("G2").Value = 80959.25712661834
$ws.Range("I2").Value = 161710.6685703679
$ws.Range("L2").Value = 484922.2142001599
$ws.Range("M2").Value = 105953.7713982
$ws.Range("N2").Value = 70003.73489578845
$ws.Range("O2").Value = 69744.89343456978

# ---------------------------------------------------------------------------
# Sheet "2030"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item(2)
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 31203.23858116339
$ws.Range("E2").Value = 170658.5511254234
$ws.Range("I2").Value = 209080.6134235085
$ws.Range("L2").Value = 63518.11613148725
$ws.Range("M2").Value = 68536.72857011756
$ws.Range("N2").Value = 19285.19160463996
$ws.Range("O2").Value = 27033.1386905727

# ---------------------------------------------------------------------------
# Sheet "2035"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item(3)
$ws.Range("A2").Value = 27543.1755456332
$ws.Range("B2").Value = 22113.21643273498
$ws.Range("E2").Value = 114655.4402706629
$ws.Range("I2").Value = 153866.0861464091
$ws.Range("M2").Value = 44638.22942194272
$ws.Range("N2").Value = 39676.88529639924
$ws.Range("O2").Value = 31311.04369977792

# ---------------------------------------------------------------------------
# Sheet "2040"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item(4)
$ws.Range("N2").Value = 1142.580190039942
$ws.Range("O2").Value = 0

# ---------------------------------------------------------------------------
# Sheet "2045"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item(5)
$ws.Range("A2").Value = 29588.33508286276
$ws.Range("N2").Value = 4347.543515635315
$ws.Range("O2").Value = 20429.76977394434

# Sheet "2050" has no changes.
